$d = $word.ActiveDocument

# Correct the "${HP}" merge-field tag to "${HPP}" in the compressor cost
# paragraph (it was a typo - should match the other HP-prefixed tags).
# Scope the Find to a narrow Range around the match so only that run is
# touched, instead of letting a whole-story Find.Execute merge every
# neighboring run in the paragraph together.
$rng = $d.Content.Duplicate
$rng.Find.Execute("`${HP}", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
if ($rng.Find.Found) {
    $rng.Text = "`${HPP}"
}
